# Apply edits to match the target workbook state:
# - Reword the question in row 2/4 and add dedicated "keyword" rows (futebol/boxe)
# - Append new Q&A rows for promotions, bingo results, recommendations, live betting,
#   payment methods, withdrawal time/keyword and talking to an agent
# - Update row heights for the reflowed long-answer rows
# - Widen the two columns and refresh the sheet view (zoom/selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- String literals used throughout the sheet ---
$s_perguntas = 'Perguntas '
$s_respostas = 'Respostas'
$s_futebol_long = 'Existem várias formas de apostar em futebol, algumas das mais comuns incluem:

Aposta 1x2: Você aposta no resultado da partida, sendo "1" para vitória do time da casa, "X" para empate e "2" para vitória do time visitante.

Aposta no número de gols (Over/Under): Aposte se haverá mais ou menos gols do que um número específico (geralmente 2,5). Por exemplo, se apostar no "over 2,5", você ganha se a partida terminar com 3 ou mais gols.

Aposta dupla chance: Você aposta em duas opções ao mesmo tempo, como "1 ou X" ou "X ou 2", aumentando suas chances de vitória.

Apostas em handicap: Neste tipo de aposta, um time começa com uma vantagem ou desvantagem fictícia para equilibrar as odds entre os times. Por exemplo, se você apostar "Time A -1", isso significa que o time A precisa vencer por mais de um gol de diferença para você ganhar a aposta.'
$s_futebol_kw = 'futebol'
$s_boxe_kw = 'boxe'
$s_boxe_long = 'Aposta no vencedor (Moneyline): A aposta mais simples e comum. Você aposta em quem vai ganhar a luta, seja o lutador A ou o lutador B.

Aposta no método de vitória: Você pode apostar em como o lutador ganhará a luta:

Knockout (KO): Se você acha que um dos lutadores vai vencer por nocaute.

Nocaute Técnico (TKO): Se o lutador vai vencer por uma interrupção do árbitro (geralmente por lesão ou incapacidade de continuar).

Decisão: Se a luta vai até o final dos rounds e um dos lutadores vencer por pontos, após decisão dos juízes.

Desqualificação ou outro resultado: Se o lutador vencer por desqualificação ou outro motivo.

Aposta no round de vitória: Se você acha que um dos lutadores vai ganhar em um determinado round, pode apostar qual round será esse.

Aposta no número de rounds (Over/Under): O número de rounds define a duração da luta. Por exemplo, se você apostar no "Over 8,5", a luta deve ir além do 8º round para que você vença a aposta. Se apostar no "Under 8,5", a luta precisa acabar antes do 9º round.

Aposta em empate: Embora muito menos comum, você pode apostar que a luta terminará em empate. Esse tipo de aposta tem odds bem altas, pois empates são raros no boxe.'
$s_rollover_q = ' O que é rollover ?'
$s_rollover_kw = 'rollover'
$s_rollover_long = 'Rollover é um termo frequentemente usado em apostas e cassinos online, mas também aparece em outros contextos financeiros. No universo das apostas, o rollover refere-se à quantidade de vezes que você precisa apostar o valor de um bônus ou depósito antes de poder retirar qualquer valor ganho. Ele funciona como uma condição que as casas de apostas ou cassinos exigem para que o bônus seja liberado e transformado em saldo disponível para saque.

Como funciona o rollover nas apostas?
Exemplo de Rollover em um Bônus: Suponha que você receba um bônus de R$100 após fazer um depósito de R$100 em uma casa de apostas. A casa pode exigir um rollover de 5x (cinco vezes) sobre o valor do bônus. Isso significa que, antes de poder retirar qualquer valor, você precisa apostar o valor do bônus (R$100) 5 vezes, ou seja, R$500 em apostas.

Por que o Rollover é importante?

Condicional: O rollover é uma maneira de garantir que o jogador esteja envolvido com as apostas antes de retirar o valor ganho no bônus. Isso evita que as pessoas façam apenas um depósito, ganhem o bônus e depois retirem sem realmente apostarem ou jogarem.

Variação: O valor do rollover pode variar dependendo da casa de apostas e do tipo de bônus. Alguns bônus podem ter rollovers mais baixos, enquanto outros têm requisitos mais altos.

Tipos de Rollover:

Bônus de depósito: Onde você recebe um bônus baseado no valor do depósito, e o rollover se aplica sobre o valor do bônus (não o depósito).

Apostas qualificadas: Algumas casas exigem que você aposte em mercados ou tipos de apostas específicas para contar no rollover.

Rollover para bônus de apostas gratuitas: Se você receber apostas grátis, também pode haver um rollover específico para que você possa usar os ganhos dessas apostas para fazer retiradas.'
$s_odds_long = 'Odds (ou probabilidades) são uma representação matemática das chances de um evento ocorrer em uma aposta. Elas indicam não só a probabilidade de um resultado específico, mas também quanto você pode ganhar em relação ao valor apostado. As odds são fundamentais para entender o risco de uma aposta e os possíveis retornos.

As odds podem ser apresentadas de três formas principais: decimais, fracionárias e americanas. Vamos explicar cada uma delas:

1. Odds Decimais (ou "Odds Européias")
As odds decimais são as mais comuns em muitas casas de apostas online e são muito fáceis de entender. Elas indicam o retorno total (incluindo o valor apostado) para cada unidade apostada.

Como funcionam:

Se você apostar R$100 em uma odd de 2.50, isso significa que, se você ganhar, você receberá R$250 de volta (R$100 x 2.50 = R$250), o que inclui o valor da aposta original.

O número 2.50 reflete que, para cada R$1 apostado, você receberá R$2,50 (ou seja, seu investimento inicial + o lucro).

Exemplo:

Aposta: R$100

Odds: 2.50

Retorno (se ganhar): R$250 (R$100 x 2.50)

2. Odds Fracionárias (ou "Odds Britânicas")
As odds fracionárias são muito usadas em apostas tradicionais, especialmente no Reino Unido. Elas mostram o lucro que você pode obter em relação à sua aposta.

Como funcionam:

Se você vê 5/1 (cinco para um), isso significa que, para cada R$1 apostado, você ganha R$5 de lucro (além de receber o valor da sua aposta inicial de volta).

Se a odd for 1/5 (um para cinco), isso significa que, para cada R$5 apostados, você ganha R$1 de lucro (novamente, o valor da aposta original será devolvido).

Exemplo:

Aposta: R$100

Odds: 5/1

Retorno (se ganhar): R$600 (R$100 x 5 + a aposta inicial de R$100)

3. Odds Americanas (ou "Moneyline")
As odds americanas podem ser positivas ou negativas, dependendo de quem é o favorito ou o azarão.

Odds positivas (+): As odds positivas mostram o quanto você ganharia em lucro para cada R$100 apostados. Por exemplo, +200 significa que, para cada R$100 apostados, você ganhará R$200 de lucro (sem contar a sua aposta original).

Odds negativas (-): As odds negativas mostram quanto você precisa apostar para ganhar R$100 de lucro. Por exemplo, -150 significa que você precisa apostar R$150 para ganhar R$100 de lucro (sem contar a sua aposta original).

Exemplo de odds positivas:

Aposta: R$100

Odds: +200

Retorno (se ganhar): R$300 (R$100 apostados + R$200 de lucro)

Exemplo de odds negativas:

Aposta: R$100

Odds: -150

Retorno (se ganhar): R$166,67 (R$100 apostados + R$66,67 de lucro)

Como as odds refletem a probabilidade?
As odds também estão intimamente ligadas à probabilidade de um evento ocorrer. Quanto menores as odds, maior a probabilidade do evento acontecer. Por outro lado, odds mais altas indicam um evento mais improvável de ocorrer.

Exemplo:

Odds de 1.20 indicam uma probabilidade alta de sucesso (aproximadamente 83% de chance).

Odds de 5.00 indicam uma probabilidade baixa de sucesso (aproximadamente 20% de chance).

Como calcular o lucro?
Para calcular o lucro de uma aposta, você pode usar a fórmula simples:

Lucro = Valor apostado x (Odds - 1) (se for odd decimal)

Por exemplo, se você apostar R$100 com odds de 3.00:

Lucro = 100 x (3.00 - 1) = 100 x 2 = R$200

Ou seja, você ganharia R$200 de lucro e o retorno total seria R$300 (R$200 de lucro + R$100 da aposta inicial).'
$s_odds_q = 'O que são odds ?'
$s_odds_kw = 'odds '
$s_futebol_rules_q = 'quais as regras do futebol?
'
$s_boxe_rules_q = 'quais as regras do boxe?
'
$s_promo_q = 'Quais sãos as promoções de hoje ?'
$s_promo_a = 'Opa você vai adorar!!Temos a ultra aposta, aposta turbi,bingo da sorte, campo minado e muito mais !'
$s_bingo_q = 'Qual o resultado de bingo de hoje?'
$s_bingo_a = 'O resultado foi : Dona Maria, José Melo,Luis Bruno,Antonio Buono e Charle Manck'
$s_recommend_q = 'Pode me recomendar no que apostar ?'
$s_recommend_a = 'Claro !Futebol,MMA,Basquete,tênis,beisebol'
$s_live_q = 'Você tem apostas live?'
$s_live_a = 'Neste momento basquete e futebol ,segue o link!'
$s_payment_q = 'Quais as formas de pagamento você aceita?'
$s_payment_a = 'Pix e depósito bancario, em 10 minutos você já consegue apostar !'
$s_withdraw_q = 'Qual o prazo para eu sacar meus ganhos ?'
$s_withdraw_a = '3 dias úteis'
$s_withdraw_kw = 'saque'
$s_agent_q = 'Quero falar com um analista ou atendente ?'
$s_agent_a = 'Claro me confirme o seu nome e e-mail '

# --- Row 1: header (unchanged) ---
$ws.Range("A1").Value = $s_perguntas
$ws.Range("B1").Value = $s_respostas

# --- Row 2 ---
$ws.Range("A2").Value = $s_futebol_rules_q
$ws.Range("B2").Value = $s_futebol_long
$ws.Range("A2").WrapText = $true
$ws.Range("B2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 88.2

# --- Row 3 ---
$ws.Range("A3").Value = $s_futebol_kw
$ws.Range("B3").Value = $s_futebol_long
$ws.Range("B3").WrapText = $true
$ws.Rows.Item(3).RowHeight = 100.2

# --- Row 4 ---
$ws.Range("A4").Value = $s_boxe_rules_q
$ws.Range("B4").Value = $s_boxe_long
$ws.Range("A4").WrapText = $true
$ws.Range("B4").WrapText = $true
$ws.Rows.Item(4).RowHeight = 345.6

# --- Row 5 ---
$ws.Range("A5").Value = $s_boxe_kw
$ws.Range("B5").Value = $s_boxe_long
$ws.Range("B5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 345.6

# --- Row 6 ---
$ws.Range("A6").Value = $s_rollover_q
$ws.Range("B6").Value = $s_rollover_long
$ws.Range("A6").WrapText = $true
$ws.Range("B6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 409.6

# --- Row 7 ---
$ws.Range("A7").Value = $s_rollover_kw
$ws.Range("B7").Value = $s_rollover_long
$ws.Range("B7").WrapText = $true
$ws.Rows.Item(7).RowHeight = 409.6

# --- Row 8 ---
$ws.Range("A8").Value = $s_odds_q
$ws.Range("B8").Value = $s_odds_long
$ws.Range("A8").WrapText = $true
$ws.Range("B8").WrapText = $true
$ws.Rows.Item(8).RowHeight = 409.6

# --- Row 9 ---
$ws.Range("A9").Value = $s_odds_kw

# --- Row 10 ---
$ws.Range("A10").Value = $s_promo_q
$ws.Range("B10").Value = $s_promo_a
$ws.Range("A10").WrapText = $true
$ws.Range("B10").WrapText = $true

# --- Row 11 ---
$ws.Range("A11").Value = $s_bingo_q
$ws.Range("B11").Value = $s_bingo_a
$ws.Range("B11").WrapText = $true

# --- Row 12 ---
$ws.Range("A12").Value = $s_recommend_q
$ws.Range("B12").Value = $s_recommend_a
$ws.Range("A12").WrapText = $true
$ws.Range("B12").WrapText = $true

# --- Row 13 ---
$ws.Range("A13").Value = $s_live_q
$ws.Range("B13").Value = $s_live_a
$ws.Range("B13").WrapText = $true

# --- Row 14 ---
$ws.Range("A14").Value = $s_payment_q
$ws.Range("B14").Value = $s_payment_a
$ws.Range("A14").WrapText = $true
$ws.Range("B14").WrapText = $true

# --- Row 15 ---
$ws.Range("A15").Value = $s_withdraw_q
$ws.Range("B15").Value = $s_withdraw_a
$ws.Range("B15").WrapText = $true

# --- Row 16 ---
$ws.Range("A16").Value = $s_withdraw_kw
$ws.Range("B16").Value = $s_withdraw_a
$ws.Range("A16").WrapText = $true
$ws.Range("B16").WrapText = $true

# --- Row 17 ---
$ws.Range("A17").Value = $s_agent_q
$ws.Range("B17").Value = $s_agent_a
$ws.Range("B17").WrapText = $true

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 45.59
$ws.Columns.Item(2).ColumnWidth = 88.25

# --- Sheet view: zoom + selection on the new last (empty) row ---
$excel.ActiveWindow.Zoom = 90
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A18").Select()
